$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two trailing test rows (old rows 4 and 5) - this also
# shrinks the sheet's used range / dimension down to A1:C3.
$ws.Range("A4:C5").Delete()

# Row 2: update the existing submission's values.
$ws.Range("A2").Value = "Test User Fixed"
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = "First submission after fix"

# Row 3: new submission replacing the old "John Doe" row. Note the
# Age column here is stored as text ("25"), not a number, so force
# text formatting while entering it, then drop back to the default
# "Normal" style so no stray number format lingers on the cell.
$ws.Range("A3").Value = "Second Test"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "25"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "Testing multiple entries"
